{"js": "// Load all paragraphs in the document body so we can locate the\n// \"Meta description\" paragraph (right after the title) and the final\n// paragraph (which currently holds the image-generation \"Prompt\" text).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- Part 1: remove the \"Meta description: ...\" paragraph that sits\n// right below the H1 title at the top of the document. ---\nconst metaIndex = paragraphs.items.findIndex((p) =>\n  p.text.trim().startsWith(\"Meta description\")\n);\nif (metaIndex !== -1) {\n  paragraphs.items[metaIndex].delete();\n  await context.sync();\n}\n\n// --- Part 2: at the end of the document, insert a new bold paragraph\n// with the page title, and replace the old \"Prompt: ...\" text of the\n// last paragraph with the meta-description copy (keeping its italic\n// formatting). ---\n\n// Re-fetch paragraphs since the collection changed after the delete above.\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs2.items;\nconst lastParagraph = items[items.length - 1];\n\n// Insert a new paragraph right before the last one, containing the bold title.\nconst titleParagraph = lastParagraph.insertParagraph(\n  \"Play Frost Queen Jackpots Free - Review 2021\",\n  Word.InsertLocation.before\n);\ntitleParagraph.font.bold = true;\ntitleParagraph.font.italic = false;\n\n// Replace the text of the last paragraph (previously the \"Prompt: ...\"\n// text) with the review meta description, keeping its existing (italic)\n// character formatting intact.\nlastParagraph.getRange(\"Whole\").insertText(\n  \"Read our review of Frost Queen Jackpots. Try this game for free and discover the numerous bonus features and potential payouts.\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Apply the \"Added a few more slots\" edit:\n#  1. Remove the \"Meta description: ...\" paragraph that sits right under\n#     the H1 title at the top of the document.\n#  2. At the end of the document, insert a new bold paragraph repeating\n#     the page title, and replace the old \"Prompt: ...\" text of the final\n#     paragraph with the meta-description copy (keeping its italic run\n#     formatting).\n\n$d = $word.ActiveDocument\n\n# --- Part 1: delete the \"Meta description\" paragraph -----------------\n$metaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $pText = $d.Paragraphs.Item($i).Range.Text\n    if ($pText.TrimStart().StartsWith(\"Meta description\")) {\n        $metaIndex = $i\n        break\n    }\n}\nif ($metaIndex -ne -1) {\n    $d.Paragraphs.Item($metaIndex).Range.Delete()\n}\n\n# --- Part 2: add the bold title paragraph before the last paragraph, --\n# --- then replace the last paragraph's text with the review blurb. ---\n$n = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($n)\n$insertionPoint = $lastPara.Range\n$insertionPoint.Collapse(1)  # wdCollapseStart\n\n$titleText = \"Play Frost Queen Jackpots Free - Review 2021\"\n$insertionPoint.InsertBefore($titleText + \"`r\")\n\n# The newly inserted paragraph is now at the old last-paragraph's index;\n# it pushed the former last paragraph (the \"Prompt: ...\" one) one slot down.\n$titlePara = $d.Paragraphs.Item($n)\n$titleRange = $titlePara.Range.Duplicate\n[void]$titleRange.MoveEnd(1, -1)  # wdCharacter; exclude the paragraph mark\n$titleRange.Font.Bold = 1\n$titleRange.Font.Italic = 0\n\n$finalPara = $d.Paragraphs.Item($n + 1)\n$finalRange = $finalPara.Range.Duplicate\n[void]$finalRange.MoveEnd(1, -1)  # wdCharacter; exclude the paragraph mark\n$finalRange.Text = \"Read our review of Frost Queen Jackpots. Try this game for free and discover the numerous bonus features and potential payouts.\"\n"}
